$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update K4 value (0.1459 -> 0.0902)
$ws.Range("K4").Value = 0.0902

# Add new "End Date" values in column N for rows 2-10 (46022 = 12/31/2025)
$ws.Range("N2").Value = 46022
$ws.Range("N3").Value = 46022
$ws.Range("N4").Value = 46022
$ws.Range("N5").Value = 46022
$ws.Range("N6").Value = 46022
$ws.Range("N7").Value = 46022
$ws.Range("N8").Value = 46022
$ws.Range("N9").Value = 46022
$ws.Range("N10").Value = 46022

# Update the active selection from M7 to L4
$ws.Range("L4").Select()
